# Updates the cryptos list prices/volumes, as published by the "Updated cryptos
# list" GitHub Actions scraper job. Only the data cells change: column B (Coin
# name), column C (coinranking.com link), column D (Price) and column E
# (Volume(1h)); column A (the 0-based rank index) is left untouched.
#
# A handful of rows (24/25, 30/31, 40/41, 45/46, 48/49) swap rank order between
# two coins, so those pairs get their B/C/D/E cells fully replaced rather than
# just the numbers refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.150.91"
$ws.Range("E2").Value = "  -7.09%  "

# Row 3
$ws.Range("D3").Value = "3.228.65"
$ws.Range("E3").Value = "  -9.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").Value = "'173.70"
$ws.Range("E5").Value = "  -16.31%  "

# Row 6
$ws.Range("D6").Value = "'506.67"
$ws.Range("E6").Value = "  -10.86%  "

# Row 7
$ws.Range("D7").Value = "'0.581"
$ws.Range("E7").Value = "  -4.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "3.218.15"
$ws.Range("E9").Value = "  -9.74%  "

# Row 10
$ws.Range("D10").Value = "'0.604"
$ws.Range("E10").Value = "  -10.70%  "

# Row 11
$ws.Range("D11").Value = "'56.20"
$ws.Range("E11").Value = "  -8.77%  "

# Row 12
$ws.Range("E12").Value = "  -13.05%  "

# Row 13
$ws.Range("E13").Value = "  -10.73%  "

# Row 14
$ws.Range("D14").Value = "'8.90"
$ws.Range("E14").Value = "  -12.73%  "

# Row 15
$ws.Range("D15").Value = "3.750.43"
$ws.Range("E15").Value = "  -9.24%  "

# Row 16
$ws.Range("E16").Value = "  -7.15%  "

# Row 17
$ws.Range("D17").Value = "3.232.03"
$ws.Range("E17").Value = "  -9.10%  "

# Row 18
$ws.Range("D18").Value = "62.955.14"
$ws.Range("E18").Value = "  -6.97%  "

# Row 19
$ws.Range("D19").Value = "'16.85"
$ws.Range("E19").Value = "  -11.56%  "

# Row 20
$ws.Range("D20").Value = "'10.56"
$ws.Range("E20").Value = "  -13.21%  "

# Row 21
$ws.Range("D21").Value = "'0.924"
$ws.Range("E21").Value = "  -12.62%  "

# Row 22
$ws.Range("D22").Value = "'362.42"
$ws.Range("E22").Value = "  -9.57%  "

# Row 23
$ws.Range("D23").Value = "'78.39"
$ws.Range("E23").Value = "  -6.74%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'3.56"
$ws.Range("E24").Value = "  -14.35%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'10.67"
$ws.Range("E25").Value = "  -14.41%  "

# Row 26
$ws.Range("D26").Value = "'3.70"
$ws.Range("E26").Value = "  -6.05%  "

# Row 27
$ws.Range("E27").Value = "  -9.92%  "

# Row 28
$ws.Range("D28").Value = "'10.99"
$ws.Range("E28").Value = "  -11.21%  "

# Row 29
$ws.Range("D29").Value = "'8.11"
$ws.Range("E29").Value = "  -11.99%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'27.81"
$ws.Range("E30").Value = "  -11.50%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'631.82"
$ws.Range("E31").Value = "  -5.68%  "

# Row 32
$ws.Range("D32").Value = "'6.50"
$ws.Range("E32").Value = "  -14.71%  "

# Row 33
$ws.Range("E33").Value = "  -9.71%  "

# Row 34
$ws.Range("D34").Value = "'58.41"
$ws.Range("E34").Value = "  -7.56%  "

# Row 35
$ws.Range("E35").Value = "  -10.12%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("D37").Value = "'35.03"
$ws.Range("E37").Value = "  -14.05%  "

# Row 38
$ws.Range("D38").Value = "'0.369"
$ws.Range("E38").Value = "  -9.55%  "

# Row 39
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.823.09"
$ws.Range("E40").Value = "  -10.63%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  -9.72%  "

# Row 42
$ws.Range("E42").Value = "  -15.12%  "

# Row 43
$ws.Range("E43").Value = "  -19.43%  "

# Row 44
$ws.Range("D44").Value = "'2.56"
$ws.Range("E44").Value = "  -8.03%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  +2.74%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.28"
$ws.Range("E46").Value = "  -14.01%  "

# Row 47
$ws.Range("D47").Value = "'0.0373"
$ws.Range("E47").Value = "  -8.67%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  -4.38%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.121"
$ws.Range("E49").Value = "  -6.45%  "

# Row 50
$ws.Range("D50").Value = "'131.12"
$ws.Range("E50").Value = "  -4.97%  "

# Row 51
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  -17.02%  "
